$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44307
$ws.Range("M2").Value = 40
$ws.Range("N2").Value = 10000
$ws.Range("O2").Value = 10000
$ws.Range("P2").Value = 10000
$ws.Range("S2").Value = 1000
$ws.Range("D3").Value = 44301
$ws.Range("D4").Value = 44333
$ws.Range("M4").Value = 58
$ws.Range("N4").Value = 10000
$ws.Range("O4").Value = 10000
$ws.Range("P4").Value = 10000
$ws.Range("S4").Value = 1000
$ws.Range("D5").Value = 44333
$ws.Range("M5").Value = 65
$ws.Range("N5").Value = 9000
$ws.Range("O5").Value = 9000
$ws.Range("P5").Value = 9000
$ws.Range("S5").Value = 900
$ws.Range("D6").Value = 44333
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 8000
$ws.Range("O6").Value = 8000
$ws.Range("P6").Value = 8000
$ws.Range("S6").Value = 800
$ws.Range("D7").Value = 44306
$ws.Range("D8").Value = 44309
$ws.Range("M8").Value = 45
$ws.Range("N8").Value = 10000
$ws.Range("O8").Value = 10000
$ws.Range("P8").Value = 10000
$ws.Range("S8").Value = 1000
$ws.Range("D9").Value = 44328
$ws.Range("M9").Value = 45
$ws.Range("N9").Value = 8000
$ws.Range("O9").Value = 8000
$ws.Range("P9").Value = 8000
$ws.Range("S9").Value = 800
$ws.Range("D10").Value = 44328
$ws.Range("M10").Value = 48
$ws.Range("N10").Value = 7000
$ws.Range("O10").Value = 7000
$ws.Range("P10").Value = 7000
$ws.Range("S10").Value = 700
$ws.Range("D11").Value = 44326
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 65
$ws.Range("D12").Value = 44326
$ws.Range("L12").Value = "Segunda"
$ws.Range("M12").Value = 67
$ws.Range("N12").Value = 8000
$ws.Range("O12").Value = 8000
$ws.Range("P12").Value = 8000
$ws.Range("S12").Value = 800
$ws.Range("D13").Value = 44319
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 68
$ws.Range("N13").Value = 10000
$ws.Range("O13").Value = 10000
$ws.Range("P13").Value = 10000
$ws.Range("S13").Value = 1000
$ws.Range("D14").Value = 44319
$ws.Range("L14").Value = "Segunda"
$ws.Range("M14").Value = 57
$ws.Range("N14").Value = 8000
$ws.Range("O14").Value = 8000
$ws.Range("P14").Value = 8000
$ws.Range("S14").Value = 800
$ws.Range("D15").Value = 44343
$ws.Range("L15").Value = "Especial"
$ws.Range("M15").Value = 47
$ws.Range("N15").Value = 10000
$ws.Range("O15").Value = 10000
$ws.Range("P15").Value = 10000
$ws.Range("R15").Value = "Región Metropolitana"
$ws.Range("S15").Value = 1000
$ws.Range("D16").Value = 44343
$ws.Range("M16").Value = 50
$ws.Range("N16").Value = 9000
$ws.Range("O16").Value = 9000
$ws.Range("P16").Value = 9000
$ws.Range("R16").Value = "Región Metropolitana"
$ws.Range("S16").Value = 900
$ws.Range("D17").Value = 44343
$ws.Range("M17").Value = 58
$ws.Range("R17").Value = "Región Metropolitana"
$ws.Range("D18").Value = 44312
$ws.Range("M18").Value = 48
$ws.Range("N18").Value = 10000
$ws.Range("O18").Value = 10000
$ws.Range("P18").Value = 10000
$ws.Range("R18").Value = "Provincia de Quillota"
$ws.Range("S18").Value = 1000
$ws.Range("D19").Value = 44314
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 47
$ws.Range("N19").Value = 9000
$ws.Range("O19").Value = 9000
$ws.Range("P19").Value = 9000
$ws.Range("R19").Value = "Provincia de Quillota"
$ws.Range("S19").Value = 900
$ws.Range("D20").Value = 44308
$ws.Range("D21").Value = 44308
$ws.Range("L21").Value = "Segunda"
$ws.Range("M21").Value = 48
$ws.Range("N21").Value = 8000
$ws.Range("O21").Value = 8000
$ws.Range("P21").Value = 8000
$ws.Range("S21").Value = 800
$ws.Range("D22").Value = 44699
$ws.Range("M22").Value = 56
$ws.Range("N22").Value = 12000
$ws.Range("O22").Value = 12000
$ws.Range("P22").Value = 12000
$ws.Range("R22").Value = "Provincia de Quillota"
$ws.Range("S22").Value = 1200
$ws.Range("D23").Value = 44699
$ws.Range("M23").Value = 60
$ws.Range("N23").Value = 10000
$ws.Range("O23").Value = 10000
$ws.Range("P23").Value = 10000
$ws.Range("R23").Value = "Provincia de Quillota"
$ws.Range("S23").Value = 1000
$ws.Range("D24").Value = 44329
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 56
$ws.Range("N24").Value = 9000
$ws.Range("O24").Value = 9000
$ws.Range("P24").Value = 9000
$ws.Range("S24").Value = 900
$ws.Range("D25").Value = 44329
$ws.Range("L25").Value = "Segunda"
$ws.Range("M25").Value = 50
$ws.Range("N25").Value = 8000
$ws.Range("O25").Value = 8000
$ws.Range("P25").Value = 8000
$ws.Range("R25").Value = "Región Metropolitana"
$ws.Range("S25").Value = 800
$ws.Range("D26").Value = 44322
$ws.Range("L26").Value = "Primera"
$ws.Range("M26").Value = 56
$ws.Range("N26").Value = 10000
$ws.Range("O26").Value = 10000
$ws.Range("P26").Value = 10000
$ws.Range("S26").Value = 1000
$ws.Range("D27").Value = 44322
$ws.Range("L27").Value = "Segunda"
$ws.Range("M27").Value = 40
$ws.Range("N27").Value = 8000
$ws.Range("O27").Value = 8000
$ws.Range("P27").Value = 8000
$ws.Range("S27").Value = 800
$ws.Range("D28").Value = 44315
$ws.Range("L28").Value = "Primera"
$ws.Range("M28").Value = 45
$ws.Range("N28").Value = 10000
$ws.Range("O28").Value = 10000
$ws.Range("P28").Value = 10000
$ws.Range("S28").Value = 1000
$ws.Range("D29").Value = 44302
$ws.Range("M29").Value = 45
$ws.Range("D30").Value = 44321
$ws.Range("M30").Value = 58
$ws.Range("N30").Value = 9000
$ws.Range("O30").Value = 9000
$ws.Range("P30").Value = 9000
$ws.Range("S30").Value = 900
$ws.Range("D31").Value = 44323
$ws.Range("L31").Value = "Primera"
$ws.Range("M31").Value = 60
$ws.Range("N31").Value = 10000
$ws.Range("O31").Value = 10000
$ws.Range("P31").Value = 10000
$ws.Range("S31").Value = 1000
$ws.Range("D32").Value = 44323
$ws.Range("L32").Value = "Segunda"
$ws.Range("M32").Value = 50
$ws.Range("N32").Value = 9000
$ws.Range("O32").Value = 9000
$ws.Range("P32").Value = 9000
$ws.Range("S32").Value = 900
